$wb = $excel.ActiveWorkbook

# --- "bike list": insert a new column before B, add a photo/source column ---
$ws3 = $wb.Worksheets.Item("bike list")
$ws3.Columns("B:B").Insert()
$ws3.Range("B5").Value = "photo"
$ws3.Range("B6").Value = "Suzuki.ca"

# --- "edit tips": drop a dev note about removing the picture/intro from the forum post ---
$ws2 = $wb.Worksheets.Item("edit tips")
$ws2.Range("J38").Value = "q.remove_pictureintro!"

# Leave the cursor/selection on "bike list" where the new column was edited ...
$ws3.Activate()
$ws3.Range("B7").Select()

# ... then return to "edit tips" (the tab that stays active) with its own
# last-used selection.
$ws2.Activate()
$ws2.Range("L33").Select()
